# Applies the "Saldo.xlsx" update:
#  - changes LAILA's balance from 350000 to 300000
#  - inserts several new account rows
#  - removes a few stale rows whose balances moved earlier in the sheet
#
# All row numbers below refer to the ORIGINAL (before-edit) layout and are
# processed from the bottom of the sheet upward so that earlier operations
# never invalidate the row numbers used by later ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DataRow($rowNum, $conta, $nome, $saldo) {
    # Force column A to text so account numbers keep their leading zeros
    # (matches the source data, which is stored as inlineStr).
    $ws.Cells.Item($rowNum, 1).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 1).Value = $conta
    $ws.Cells.Item($rowNum, 2).Value = $nome
    $ws.Cells.Item($rowNum, 3).Value = $saldo
}

# --- Deletions (bottom to top) ---------------------------------------

# Row 491: 004381194 ALINNE 0.12  (stale tail entry; account now appears
# earlier in the sheet with an updated balance - see insertion below)
$ws.Rows.Item(491).Delete()

# Row 278: 004724018 ASPA 29.83  (stale tail entry; account now appears
# earlier in the sheet with an updated balance - see insertion below)
$ws.Rows.Item(278).Delete()

# Row 14: 004504449 KELMA 1013.75  (removed entirely)
$ws.Rows.Item(14).Delete()

# --- Insertions (bottom to top) ---------------------------------------

# Two new rows before row 8 (004467884 ANA). The final desired order is
# ANDRE (row 8) then JOSE (row 9), so insert/fill JOSE first - the next
# insert above it will push it down to row 9, then fill ANDRE on top at
# row 8.
$ws.Rows.Item(8).Insert()
Set-DataRow 8 "004639776" "JOSE" 12000

$ws.Rows.Item(8).Insert()
Set-DataRow 8 "005639781" "ANDRE" 13000

# New row before row 7 (005274028 RAFAEL)
$ws.Rows.Item(7).Insert()
Set-DataRow 7 "005135105" "BRENNER" 14626.24

# New row before row 6 (004364200 BLOCO)
$ws.Rows.Item(6).Insert()
Set-DataRow 6 "004724018" "ASPA" 22820.44

# New row before row 5 (005064129 THIAGO)
$ws.Rows.Item(5).Insert()
Set-DataRow 5 "004381194" "ALINNE" 43218.5

# New row before row 4 (004550750 THEO)
$ws.Rows.Item(4).Insert()
Set-DataRow 4 "004361159" "HFR" 110816.36

# --- Simple value edit ---------------------------------------------

# Row 2: 004641487 LAILA balance 350000 -> 300000
$ws.Cells.Item(2, 3).Value = 300000
